$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Update shared string content (order matters for shared-string table layout) ---
# 1) New destination text in sheet2 (Journey Planner)
$ws2.Range("D3").Value = "JURONG EAST"

# 2) Renamed test case ids on sheet1 (Search)
$ws1.Range("A3").Value = "verifyPostalCodeSearch_TC_002"
$ws1.Range("A2").Value = "verifyKeywordSearch_TC_001"

# 3) Renamed test case ids on sheet2 (Journey Planner)
$ws2.Range("A2").Value = "verifyRouteToDestination1_TC_003"
$ws2.Range("A3").Value = "verifyRouteToDestination2_TC_004"

# 4) New test row on sheet2
$ws2.Range("A4").Value = "verifySwitchDestination_TC_005"
$ws2.Range("B4").Value = "'339780"
$ws2.Range("C4").Value = "JURONG EAST"

# --- Column width tweak on sheet2 ---
$ws2.Columns.Item(1).ColumnWidth = 36.7

# --- Selections: leave sheet1's cursor on B8, and sheet2 (the active tab) on A5 ---
$ws1.Range("B8").Select() | Out-Null
$ws2.Range("A5").Select() | Out-Null
